$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.567.29"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "3.344.94"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'256.44"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'646.08"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("D7").Value = "'1.59"
$ws.Range("E7").Value = "  +15.80%  "
$ws.Range("D8").Value = "'0.408"
$ws.Range("E8").Value = "  +6.09%  "
$ws.Range("D9").Value = "'1.11"
$ws.Range("E9").Value = "  +26.26%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "3.341.65"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "'44.37"
$ws.Range("E12").Value = "  +24.08%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.208"
$ws.Range("E13").Value = "  +4.45%  "
$ws.Range("D14").Value = "98.394.36"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "'0.0000252"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").Value = "3.967.23"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "'5.56"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "3.338.96"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'7.13"
$ws.Range("E19").Value = "  +15.65%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'16.92"
$ws.Range("E20").Value = "  +12.53%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'554.74"
$ws.Range("E21").Value = "  +13.11%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").Value = "'3.55"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'10.17"
$ws.Range("E23").Value = "  +9.04%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "'0.459"
$ws.Range("E24").Value = "  +62.48%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000202"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").Value = "'6.28"
$ws.Range("E26").Value = "  +6.82%  "
$ws.Range("D27").Value = "'100.77"
$ws.Range("E27").Value = "  +12.79%  "
$ws.Range("D28").Value = "'12.63"
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.168"
$ws.Range("E29").Value = "  +36.95%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "3.516.35"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "'11.37"
$ws.Range("E31").Value = "  +21.32%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'0.190"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "'30.10"
$ws.Range("E35").Value = "  +8.52%  "
$ws.Range("D36").Value = "'0.532"
$ws.Range("E36").Value = "  +15.28%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").Value = "'2.12"
$ws.Range("E37").Value = "  +6.67%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'7.73"
$ws.Range("E38").Value = "  +5.24%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.155"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").Value = "'528.04"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("D41").Value = "'24.71"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'1.32"
$ws.Range("E42").Value = "  +4.34%  "
$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").Value = "'3.86"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").Value = "'0.837"
$ws.Range("E44").Value = "  +8.88%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'3.20"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0387"
$ws.Range("E47").Value = "  +21.71%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'8.03"
$ws.Range("E48").Value = "  +21.81%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.03"
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'5.03"
$ws.Range("E50").Value = "  +8.90%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'162.86"
$ws.Range("E51").Value = "  +0.96%  "
